$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the vendor rows: "YP.ca", "Touch Local", "YahooUK", "Mappy",
# "PagesJaunes", "Scoot", "Yell" (rows 13-19). Deleting the rows shifts
# "Zomato" (formerly row 20) up to become the new row 13.
$ws.Range("A13:A19").EntireRow.Delete()

# Update the active selection to match the new state of the sheet.
$ws.Range("J14").Select() | Out-Null
